# Update grain orientation logic to use color codes
# Grain rules by color:
#  - HS00, HS99, HS98, HS97 -> 'mixed' (can rotate 90 degrees)
#  - HS01, HS02, HS03       -> 'fixed' (cannot rotate)
#  - Other colors default to 'fixed'
#
# The generator was re-run with this new rule, which also reshuffled /
# regenerated several rows and dropped the last two rows (20, 21).
# Rewrite rows 2-19 with the new data and clear what used to be rows 20-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("柜体侧板（R）", "CB(R)-HS00-2434-574-16", 2434, 574, 16, "HS00", 1, "mixed"),
    @("柜体侧板（L）", "CB(L)-HS00-2434-574-16", 2434, 574, 16, "HS00", 1, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-704-554-16", 704, 554, 16, "HS00", 1, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-896-554-16", 896, 554, 16, "HS00", 1, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-768-554-16", 768, 554, 16, "HS00", 1, "mixed"),
    @("顶底板（单门柜体）", "DD-HS00-336-554-16", 336, 554, 16, "HS00", 2, "mixed"),
    @("层隔板（双门柜体）", "CG-HS00-896-554-25", 896, 554, 25, "HS00", 2, "mixed"),
    @("层隔板（单门柜体）", "CG-HS00-432-554-25", 432, 554, 25, "HS00", 3, "mixed"),
    @("层隔板（单门柜体）", "CG-HS00-400-554-25", 400, 554, 25, "HS00", 3, "mixed"),
    @("层隔板（双门柜体）", "CG-HS00-832-554-25", 832, 554, 25, "HS00", 3, "mixed"),
    @("层隔板（双门柜体）", "CG-HS00-736-554-25", 736, 554, 25, "HS00", 2, "mixed"),
    @("抽屉层隔板（E/F类柜体）与`"抽屉组件`"搭配使用", "CTCG-HS00-800-554-25", 800, 554, 25, "HS00", 3, "mixed"),
    @("后背板（单门柜体）", "HB-HS00-2352-336-12", 2352, 336, 12, "HS00", 1, "mixed"),
    @("门板（L/R）", "MB(R)-(门板花色)-2320-397-16", 2320, 397, 16, "(门板花色)", 2, "fixed"),
    @("门板（L/R）", "MB(R)-(门板花色)-2288-429-16", 2288, 429, 16, "(门板花色)", 2, "fixed"),
    @("门板（L/R）", "MB(L)-(门板花色)-2320-445-16", 2320, 445, 16, "(门板花色)", 1, "fixed"),
    @("底支撑(双门柜体）", "DC-HS02-704-82-16", 704, 82, 16, "HS02", 2, "fixed"),
    @("底支撑(单门柜体）", "DC-HS99-432-82-16", 432, 82, 16, "HS99", 1, "mixed")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
    $ws.Cells.Item($row, 8).Value = $item[7]
    $row = $row + 1
}

# The new data set is two rows shorter than before (rows 20 and 21 removed).
$ws.Range("A20:H21").Clear()
